# Automatic update of files.
# Removes the second "test" argument from the HYPERLINK() formulas in
# columns S, T, V, W, X, Y for rows 2-5, leaving just the URL argument.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column -> (subfolder, file extension) used to build the link target.
$linkCols = @(
    @{ Col = "S"; Folder = "artfynd";         Ext = "xlsx" },
    @{ Col = "T"; Folder = "kartor";          Ext = "png"  },
    @{ Col = "V"; Folder = "klagomål";        Ext = "docx" },
    @{ Col = "W"; Folder = "klagomålsmail";   Ext = "docx" },
    @{ Col = "X"; Folder = "tillsyn";         Ext = "docx" },
    @{ Col = "Y"; Folder = "tillsynsmail";    Ext = "docx" }
)

for ($row = 2; $row -le 5; $row++) {
    $beteckning = $ws.Range("A$row").Value2

    foreach ($entry in $linkCols) {
        $col = $entry.Col
        $folder = $entry.Folder
        $ext = $entry.Ext

        $url = "https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/$folder/$beteckning.$ext"
        $ws.Range("$col$row").Formula = "=HYPERLINK(`"$url`")"
    }
}
